# actualizado ficheiro da proposta inicial
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 2): center the header labels (was left-aligned) ---
$headerRange = $ws.Range("B2:D2")
$headerRange.HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter
$headerRange.VerticalAlignment = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignCenter

# --- Row 11: taller row, vertically centred text ---
$ws.Rows.Item(11).RowHeight = 61.8

$b11 = $ws.Range("B11")
$b11.VerticalAlignment = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignCenter

$d11 = $ws.Range("D11")
$d11.VerticalAlignment = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignCenter
$d11.WrapText = $true

# --- D16 loses its "MS" text; D23 gains it instead ---
$d16 = $ws.Range("D16")
$d16.ClearContents()
$d16.HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignGeneral
$d16.VerticalAlignment = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignBottom
$d16.WrapText = $false
$d16.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeLeft).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlLineStyleNone
$d16.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeTop).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlLineStyleNone
$d16.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeBottom).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlLineStyleNone

$d23 = $ws.Range("D23")
$d23.Value2 = "MS"

# --- Rows 18-22 in column C become one merged, centred "Exams" cell ---
$c18 = $ws.Range("C18")
$c18.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeBottom).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlLineStyleNone

$c19 = $ws.Range("C19")
$c19.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeTop).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlLineStyleNone
$c19.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeBottom).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlLineStyleNone

$c20 = $ws.Range("C20")
$c20.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeTop).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlLineStyleNone
$c20.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeBottom).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlLineStyleNone

$c21 = $ws.Range("C21")
$c21.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeTop).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlLineStyleNone
$c21.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeBottom).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlLineStyleNone

$c22 = $ws.Range("C22")
$c22.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeTop).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlLineStyleNone

$examsRange = $ws.Range("C18:C22")
$examsRange.HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter
$examsRange.VerticalAlignment = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignCenter
$examsRange.WrapText = $true
$examsRange.Merge()

# --- Page setup: landscape instead of portrait ---
$ws.PageSetup.Orientation = [Microsoft.Office.Interop.Excel.XlPageOrientation]::xlLandscape

# --- Scroll the view so row 4 is at the top when the sheet is opened ---
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
